$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("LoginPageTest")
$loginIndex = $loginSheet.Index

# Duplicate LoginPageTest -> placed right after LoginPageTest (before AdminPageTest), then rename
$loginSheet.Copy($null, $loginSheet)
$pfLogin = $wb.Worksheets.Item($loginIndex + 1)
$pfLogin.Name = "PFLoginPageTest"

# Re-fetch AdminPageTest reference (earlier references can shift after sheet insert)
$adminSheet = $wb.Worksheets.Item("AdminPageTest")

# Duplicate AdminPageTest -> placed at the very end, then rename
$adminSheet.Copy($null, $adminSheet)
$pfAdmin = $wb.Worksheets.Item($wb.Worksheets.Count)
$pfAdmin.Name = "PFAdminPageTest"

# Update the TestSuite sheet with new rows referencing the new test sheets
$testSuite = $wb.Worksheets.Item("TestSuite")

# Copy formatting from the existing data row onto the two new rows
$testSuite.Range("A3:B3").Copy()
$testSuite.Range("A4:B5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$testSuite.Range("A4").Value = "PFLoginPageTest"
$testSuite.Range("B4").Value = "Y"
$testSuite.Range("A5").Value = "PFAdminPageTest"
$testSuite.Range("B5").Value = "Y"

$testSuite.Range("A3").Select()
$pfAdmin.Select()
